# Weekly update: a new "Poroto verde" price observation is inserted as the
# first record of the Ñuble/Terminal Hortofrutícola Agro Chillán block
# (row 39), pushing every existing observation below it down by one row
# (old row 39 -> 40, ..., old row 69 -> 70).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 39; Excel shifts rows 39:69 down to 40:70.
$ws.Rows.Item(39).Insert()

# Populate the newly inserted row 39 with the new observation.
$ws.Cells.Item(39, 1).Value = 7
$ws.Cells.Item(39, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(39, 3).Value = "Ñuble"
$ws.Cells.Item(39, 4).Value = 44606
$ws.Cells.Item(39, 5).Value = 16
$ws.Cells.Item(39, 6).Value = 100112031
$ws.Cells.Item(39, 7).Value = "Poroto verde"
$ws.Cells.Item(39, 8).Value = "Sin especificar"
$ws.Cells.Item(39, 9).Value = "Primera"
$ws.Cells.Item(39, 10).Value = 60
$ws.Cells.Item(39, 11).Value = 32000
$ws.Cells.Item(39, 12).Value = 33000
$ws.Cells.Item(39, 13).Value = 32500
$ws.Cells.Item(39, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(39, 15).Value = "Región del Maule"
$ws.Cells.Item(39, 16).Value = 1300
$ws.Cells.Item(39, 17).Value = 25
$ws.Cells.Item(39, 18).Value = "Hortaliza"
